# 2007 Monthly Time Chart Drug Arrests edit
#
# Adds a new worksheet "monthly-time-chart" after the existing
# "drug-arrests-by-age-and-type-of" sheet, containing a small two-row
# table: month names across the header row and the corresponding total
# drug-arrest counts for 2007 underneath.

$wb = $excel.ActiveWorkbook

# The existing (only) worksheet stays first; the new sheet is inserted
# right after it so the workbook ends up with both tabs, in that order.
$firstSheet = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($null, $firstSheet)
$ws.Name = "monthly-time-chart"

# Row 1: "Month" label followed by the twelve month names.
$headerValues = @(
    "Month", "January", "February", "March", "April", "May", "June",
    "July", "August", "September", "October", "November", "December"
)
for ($i = 0; $i -lt $headerValues.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headerValues[$i]
}

# Row 2: the row label followed by the monthly arrest totals for 2007.
$dataValues = @(
    "Number of Drug Arrests",
    3059, 2742, 3370, 3098, 3113, 3053, 2913, 3057, 2903, 2989, 2711, 2512
)
for ($i = 0; $i -lt $dataValues.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $dataValues[$i]
}

# Match the look of the existing table: bordered text cells for the
# header row and bordered cells for the data row (re-using the
# existing/original sheet's cell formatting instead of inventing new
# style entries).
$firstSheet.Range("A1:C1").Copy()
$ws.Range("A1:M1").PasteSpecial(-4122)

$firstSheet.Range("B2:D2").Copy()
$ws.Range("A2:M2").PasteSpecial(-4122)

$ws.Columns.Item(1).ColumnWidth = 35.71

# Leave the original sheet selected/active, as it was before the edit.
$firstSheet.Activate()
